$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the upcoming period labels in column C (rows 35-36): the old "第85期"
# (7th-gen mount) entry is dropped, and subsequent periods are renumbered down
# by one (85->84+1 shift): 86->85, 87->86.
$ws.Range("C35").Value = "第85期 第四代寵物"
$ws.Range("C36").Value = "第86期 十轉技能"

# The last row's period entry (old "第87期 十轉技能") is removed entirely,
# since it is now represented by C36 above.
$ws.Range("C37").ClearContents()

# Match the saved selection/active cell state.
$ws.Range("C37").Select()

